# Översikt ÖSTRA GÖINGE.xlsx update
# 1) Bump the "Förändrad" (C column) timestamp for every data row (2..486)
#    from 2026-02-22 (46075) to 2026-02-23 (46076).
# 2) A new logging notice "A 3578-2026" (previously row 4) gets promoted to
#    row 3 with refreshed figures (more signal species / red-listed species
#    found, plus new complaint & bird-survey hyperlinks), and the notice that
#    used to sit in row 3 ("A 59155-2022") is pushed down to row 4 unchanged
#    apart from the timestamp bump.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- 1) Bump column C for all data rows ---------------------------------
$ws.Range("C2:C486").Value = 46076

# --- 2) Row 3 becomes the "A 3578-2026" notice --------------------------
$ws.Range("A3").Value = "A 3578-2026"
$ws.Range("B3").Value = 46042.71582175926
$ws.Range("D3").Value = "SKÅNE LÄN"
$ws.Range("E3").Value = "ÖSTRA GÖINGE"
$ws.Range("G3").Value = 5.9
$ws.Range("H3").Value = 1
$ws.Range("I3").Value = 4
$ws.Range("J3").Value = 0
$ws.Range("K3").Value = 1
$ws.Range("L3").Value = 0
$ws.Range("M3").Value = 0
$ws.Range("N3").Value = 0
$ws.Range("O3").Value = 1
$ws.Range("P3").Value = 1
$ws.Range("Q3").Value = 6
$ws.Range("R3").Value = "Småvänderot`r`nBlomkålssvamp`r`nBlåmossa`r`nRödgul trumpetsvamp`r`nStor revmossa`r`nKungsfågel"

$ws.Range("S3").Formula = '=HYPERLINK("https://klasma.github.io/Logging_1256/artfynd/A 3578-2026 artfynd.xlsx", "A 3578-2026")'
$ws.Range("T3").Formula = '=HYPERLINK("https://klasma.github.io/Logging_1256/kartor/A 3578-2026 karta.png", "A 3578-2026")'
$ws.Range("V3").Formula = '=HYPERLINK("https://klasma.github.io/Logging_1256/klagomål/A 3578-2026 FSC-klagomål.docx", "A 3578-2026")'
$ws.Range("W3").Formula = '=HYPERLINK("https://klasma.github.io/Logging_1256/klagomålsmail/A 3578-2026 FSC-klagomål mail.docx", "A 3578-2026")'
$ws.Range("X3").Formula = '=HYPERLINK("https://klasma.github.io/Logging_1256/tillsyn/A 3578-2026 tillsynsbegäran.docx", "A 3578-2026")'
$ws.Range("Y3").Formula = '=HYPERLINK("https://klasma.github.io/Logging_1256/tillsynsmail/A 3578-2026 tillsynsbegäran mail.docx", "A 3578-2026")'
$ws.Range("Z3").Formula = '=HYPERLINK("https://klasma.github.io/Logging_1256/fåglar/A 3578-2026 prioriterade fågelarter.docx", "A 3578-2026")'

# --- 3) Row 4 becomes the "A 59155-2022" notice -------------------------
$ws.Range("A4").Value = "A 59155-2022"
$ws.Range("B4").Value = 44897
$ws.Range("D4").Value = "SKÅNE LÄN"
$ws.Range("E4").Value = "ÖSTRA GÖINGE"
$ws.Range("G4").Value = 5.2
$ws.Range("H4").Value = 0
$ws.Range("I4").Value = 2
$ws.Range("J4").Value = 0
$ws.Range("K4").Value = 1
$ws.Range("L4").Value = 1
$ws.Range("M4").Value = 0
$ws.Range("N4").Value = 0
$ws.Range("O4").Value = 2
$ws.Range("P4").Value = 2
$ws.Range("Q4").Value = 4
$ws.Range("R4").Value = "Ask`r`nBlek kraterlav`r`nGulfotshätta`r`nLönnlav"

$ws.Range("S4").Formula = '=HYPERLINK("https://klasma.github.io/Logging_1256/artfynd/A 59155-2022 artfynd.xlsx", "A 59155-2022")'
$ws.Range("T4").Formula = '=HYPERLINK("https://klasma.github.io/Logging_1256/kartor/A 59155-2022 karta.png", "A 59155-2022")'
$ws.Range("V4").Formula = '=HYPERLINK("https://klasma.github.io/Logging_1256/klagomål/A 59155-2022 FSC-klagomål.docx", "A 59155-2022")'
$ws.Range("W4").Formula = '=HYPERLINK("https://klasma.github.io/Logging_1256/klagomålsmail/A 59155-2022 FSC-klagomål mail.docx", "A 59155-2022")'
$ws.Range("X4").Formula = '=HYPERLINK("https://klasma.github.io/Logging_1256/tillsyn/A 59155-2022 tillsynsbegäran.docx", "A 59155-2022")'
$ws.Range("Y4").Formula = '=HYPERLINK("https://klasma.github.io/Logging_1256/tillsynsmail/A 59155-2022 tillsynsbegäran mail.docx", "A 59155-2022")'

# Row height stays the sheet's fixed 15pt (the wrapped species list would
# otherwise auto-grow the row) so rows keep matching the rest of the sheet.
$ws.Rows.Item(3).RowHeight = 15
$ws.Rows.Item(4).RowHeight = 15
